$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "alias" column header (D1), styled like the other header cells
$ws.Range("D1").Value = "alias"
$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Existing rows gain an alias value in column D (no special style, like C2/C3)
$ws.Range("D2").Value = "Administrador"
$ws.Range("D3").Value = "Maicita"

# New user row (módulo de Clientes, Usuarios: iaguilera)
$ws.Range("A4").Value = "iaguilera"
$ws.Range("B4").Value = "`$2b`$10`$NKSgcAXRuxSqKkECocB2/egUEGcVN02XkKrVpbmFo3js.gGmrrtQu"
$ws.Range("C4").Value = "Ignacio Aguilera"
$ws.Range("D4").Value = "Chamo"

# Match the author's final selection in the saved file
$ws.Range("D6").Select() | Out-Null
